# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 680.7692
$ws.Range("J17").Value = 720.8333
$ws.Range("L17").Value = 2162.4999
$ws.Range("N17").Value = -2498.4999

$ws.Range("H112").Value = 720217.9399999999
$ws.Range("J112").Value = 720217.9399999999
$ws.Range("L112").Value = 2160653.82
$ws.Range("N112").Value = -2162869.82

$ws.Range("H129").Value = 387923.72
$ws.Range("I129").Value = 3020.25
$ws.Range("J129").Value = 717841
$ws.Range("K129").Value = 9060.75
$ws.Range("L129").Value = 2153523
$ws.Range("M129").Value = -4060.75
$ws.Range("N129").Value = -2163523

$ws.Range("H132").Value = 1884629.5
$ws.Range("I132").Value = 1033.4423
$ws.Range("J132").Value = 15877057
$ws.Range("K132").Value = 3100.3269
$ws.Range("L132").Value = 47631171
$ws.Range("M132").Value = -570.3269
$ws.Range("N132").Value = -47636231

$ws.Range("H135").Value = 8474999
$ws.Range("I135").Value = 367.42856
$ws.Range("J135").Value = 50000692
$ws.Range("K135").Value = 3306.85704
$ws.Range("L135").Value = 450006228
$ws.Range("M135").Value = -771.8570399999999
$ws.Range("N135").Value = -450011298

$ws.Range("H137").Value = 10602369
$ws.Range("I137").Value = 798.14
$ws.Range("J137").Value = 69499980
$ws.Range("K137").Value = 2394.42
$ws.Range("L137").Value = 208499940
$ws.Range("M137").Value = 155.5799999999999
$ws.Range("N137").Value = -208505040

$ws.Range("H138").Value = 2654.5151
$ws.Range("I138").Value = 2132.125
$ws.Range("J138").Value = 3146.1765
$ws.Range("K138").Value = 6396.375
$ws.Range("L138").Value = 9438.529500000001
$ws.Range("M138").Value = -1256.375
$ws.Range("N138").Value = -19718.5295

$ws.Range("H141").Value = 1272.6222
$ws.Range("I141").Value = 988.86487
$ws.Range("J141").Value = 2585
$ws.Range("K141").Value = 2966.59461
$ws.Range("L141").Value = 7755
$ws.Range("M141").Value = 2213.40539
$ws.Range("N141").Value = -18115

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16087.435
$ws.Range("I2").Value = 16791.408
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 16791.408
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -16678.408
$ws.Range("N2").Value = -826

$ws.Range("H37").Value = 8278.416999999999
$ws.Range("J37").Value = 12088.2
$ws.Range("L37").Value = 12088.2
$ws.Range("N37").Value = -12634.2

$ws.Range("H61").Value = 2716833.2
$ws.Range("I61").Value = 1191294.4
$ws.Range("J61").Value = 29413764
$ws.Range("K61").Value = 1191294.4
$ws.Range("L61").Value = 29413764
$ws.Range("M61").Value = -1191082.4
$ws.Range("N61").Value = -29414188

$ws.Range("H74").Value = 28396278
$ws.Range("I74").Value = 23810150
$ws.Range("J74").Value = 44447720
$ws.Range("K74").Value = 23810150
$ws.Range("L74").Value = 44447720
$ws.Range("M74").Value = -23809276
$ws.Range("N74").Value = -44449468

$ws.Range("H77").Value = 28396278
$ws.Range("I77").Value = 23810150
$ws.Range("J77").Value = 44447720
$ws.Range("K77").Value = 119050750
$ws.Range("L77").Value = 222238600
$ws.Range("M77").Value = -119046382
$ws.Range("N77").Value = -222247336

$ws.Range("H102").Value = 2105.25
$ws.Range("I102").Value = 2082.238
$ws.Range("J102").Value = 2266.3333
$ws.Range("K102").Value = 2082.238
$ws.Range("L102").Value = 2266.3333
$ws.Range("M102").Value = -460.2379999999998
$ws.Range("N102").Value = -5510.3333

$ws.Range("H116").Value = 16087.435
$ws.Range("I116").Value = 16791.408
$ws.Range("J116").Value = 600
$ws.Range("K116").Value = 16791.408
$ws.Range("L116").Value = 600
$ws.Range("M116").Value = -14497.408
$ws.Range("N116").Value = -5188

$ws.Range("H122").Value = 2641.4285
$ws.Range("I122").Value = 2248.3333
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6744.999899999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4294.999899999999
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 10686450
$ws.Range("I132").Value = 12198344
$ws.Range("J132").Value = 5051210.5
$ws.Range("K132").Value = 36595032
$ws.Range("L132").Value = 15153631.5
$ws.Range("M132").Value = -36592502
$ws.Range("N132").Value = -15158691.5

$ws.Range("H136").Value = 2716833.2
$ws.Range("I136").Value = 1191294.4
$ws.Range("J136").Value = 29413764
$ws.Range("K136").Value = 3573883.2
$ws.Range("L136").Value = 88241292
$ws.Range("M136").Value = -3571333.2
$ws.Range("N136").Value = -88246392

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16087.435
$ws.Range("I3").Value = 16791.408
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 16791.408
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = -16677.408
$ws.Range("N3").Value = -828

$ws.Range("H86").Value = 1983.94
$ws.Range("I86").Value = 1986.7959
$ws.Range("J86").Value = 1844
$ws.Range("K86").Value = 1986.7959
$ws.Range("L86").Value = 1844
$ws.Range("M86").Value = -863.7959000000001
$ws.Range("N86").Value = -4090

$ws.Range("H89").Value = 1983.94
$ws.Range("I89").Value = 1986.7959
$ws.Range("J89").Value = 1844
$ws.Range("K89").Value = 9933.979500000001
$ws.Range("L89").Value = 9220
$ws.Range("M89").Value = -4317.979500000001
$ws.Range("N89").Value = -20452

$ws.Range("H94").Value = 1348.9166
$ws.Range("I94").Value = 1069.0588
$ws.Range("J94").Value = 2028.5714
$ws.Range("K94").Value = 1069.0588
$ws.Range("L94").Value = 2028.5714
$ws.Range("M94").Value = -618.0588
$ws.Range("N94").Value = -2930.5714

$ws.Range("H134").Value = 9740951
$ws.Range("I134").Value = 9804604
$ws.Range("J134").Value = 8929371
$ws.Range("K134").Value = 29413812
$ws.Range("L134").Value = 26788113
$ws.Range("M134").Value = -29411277
$ws.Range("N134").Value = -26793183

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 405.05264
$ws.Range("I22").Value = 171.14285
$ws.Range("J22").Value = 1060
$ws.Range("K22").Value = 171.14285
$ws.Range("L22").Value = 1060
$ws.Range("M22").Value = 178.85715
$ws.Range("N22").Value = -1760

$ws.Range("H31").Value = 1361613.6
$ws.Range("I31").Value = 1091.7142
$ws.Range("J31").Value = 5690546.5
$ws.Range("K31").Value = 1091.7142
$ws.Range("L31").Value = 5690546.5
$ws.Range("M31").Value = -796.7141999999999
$ws.Range("N31").Value = -5691136.5

$ws.Range("H34").Value = 1361613.6
$ws.Range("I34").Value = 1091.7142
$ws.Range("J34").Value = 5690546.5
$ws.Range("K34").Value = 1091.7142
$ws.Range("L34").Value = 5690546.5
$ws.Range("M34").Value = -889.7141999999999
$ws.Range("N34").Value = -5690950.5

$ws.Range("H58").Value = 814036.3
$ws.Range("I58").Value = 2886.6428
$ws.Range("J58").Value = 3247485.2
$ws.Range("K58").Value = 2886.6428
$ws.Range("L58").Value = 3247485.2
$ws.Range("M58").Value = -2683.6428
$ws.Range("N58").Value = -3247891.2

$ws.Range("H107").Value = 608.525
$ws.Range("I107").Value = 244.41667
$ws.Range("J107").Value = 764.5714
$ws.Range("K107").Value = 244.41667
$ws.Range("L107").Value = 764.5714
$ws.Range("M107").Value = 1675.58333
$ws.Range("N107").Value = -4604.5714

$ws.Range("H132").Value = 623.9844000000001
$ws.Range("I132").Value = 547.2406999999999
$ws.Range("J132").Value = 1038.4
$ws.Range("K132").Value = 1641.7221
$ws.Range("L132").Value = 3115.2
$ws.Range("M132").Value = 888.2779
$ws.Range("N132").Value = -8175.200000000001

$ws.Range("H134").Value = 625764
$ws.Range("I134").Value = 778.9091
$ws.Range("J134").Value = 4445117.5
$ws.Range("K134").Value = 2336.7273
$ws.Range("L134").Value = 13335352.5
$ws.Range("M134").Value = 198.2727
$ws.Range("N134").Value = -13340422.5

$ws.Range("H136").Value = 814036.3
$ws.Range("I136").Value = 2886.6428
$ws.Range("J136").Value = 3247485.2
$ws.Range("K136").Value = 8659.928400000001
$ws.Range("L136").Value = 9742455.600000001
$ws.Range("M136").Value = -6109.928400000001
$ws.Range("N136").Value = -9747555.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2413.8
$ws.Range("I134").Value = 1684.0834
$ws.Range("J134").Value = 5332.6665
$ws.Range("K134").Value = 5052.2502
$ws.Range("L134").Value = 15997.9995
$ws.Range("M134").Value = 17.7497999999996
$ws.Range("N134").Value = -26137.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 26078
$ws.Range("I126").Value = 26078
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 78234
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -75764
$ws.Range("N126").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3090.6924
$ws.Range("I40").Value = 3311
$ws.Range("J40").Value = 3311
$ws.Range("K40").Value = 3311
$ws.Range("L40").Value = 2595
$ws.Range("M40").Value = -3175
$ws.Range("N40").Value = -2867

$ws.Range("H122").Value = 9221821
$ws.Range("I122").Value = 1184110.5
$ws.Range("J122").Value = 33334952
$ws.Range("K122").Value = 3552331.5
$ws.Range("L122").Value = 100004856
$ws.Range("M122").Value = -3549881.5
$ws.Range("N122").Value = -100009756

$ws.Range("H136").Value = 1916549.2
$ws.Range("I136").Value = 2021066.5
$ws.Range("J136").Value = 398.66666
$ws.Range("K136").Value = 6063199.5
$ws.Range("L136").Value = 1195.99998
$ws.Range("M136").Value = -6060649.5
$ws.Range("N136").Value = -6295.999980000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2980157.5
$ws.Range("I132").Value = 4274.3076
$ws.Range("J132").Value = 15875651
$ws.Range("K132").Value = 12822.9228
$ws.Range("L132").Value = 47626953
$ws.Range("M132").Value = -10292.9228
$ws.Range("N132").Value = -47632013

$ws.Range("H136").Value = 1254.5
$ws.Range("I136").Value = 441.8421
$ws.Range("J136").Value = 2219.5312
$ws.Range("K136").Value = 1325.5263
$ws.Range("L136").Value = 6658.5936
$ws.Range("M136").Value = 1224.4737
$ws.Range("N136").Value = -11758.5936

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
